$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each target cell, force text (string) number format so Excel does not
# auto-convert the numeric-looking / percent-looking text into a real number
# or percentage value, then restore the default "Normal" style so no stray
# style index is left on the cell (matches original formatting exactly).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "304.23"
Set-TextValue "E2" "0.00%"
Set-TextValue "D3" "37.11"
Set-TextValue "E3" "6.41%"
Set-TextValue "D4" "5.000"
Set-TextValue "E4" "-3.38%"
Set-TextValue "D5" "0.07864"
Set-TextValue "E5" "0.87%"
Set-TextValue "D6" "2.220"
Set-TextValue "E6" "-3.26%"
Set-TextValue "D7" "8.008"
Set-TextValue "E7" "-0.56%"
Set-TextValue "D8" "4.015"
Set-TextValue "E8" "0.60%"
Set-TextValue "D9" "0.9203"
Set-TextValue "E9" "-0.48%"
Set-TextValue "D10" "0.09604"
Set-TextValue "E10" "-4.19%"
Set-TextValue "E11" "2.90%"
Set-TextValue "D12" "0.08624"
Set-TextValue "E12" "0.53%"
Set-TextValue "D13" "0.03657"
Set-TextValue "E13" "7.68%"
Set-TextValue "D14" "0.09975"
Set-TextValue "E14" "0.61%"
Set-TextValue "D15" "0.001488"
Set-TextValue "E15" "0.29%"
Set-TextValue "D16" "0.005709"
Set-TextValue "E16" "-1.91%"
Set-TextValue "D17" "3.465"
Set-TextValue "E17" "-0.23%"
Set-TextValue "E18" "6.02%"
Set-TextValue "D19" "0.3413"
Set-TextValue "E19" "-0.08%"
Set-TextValue "E20" "-0.78%"
Set-TextValue "D21" "4.762"
Set-TextValue "E21" "4.61%"
Set-TextValue "D22" "0.2201"
Set-TextValue "E22" "-8.14%"
Set-TextValue "D23" "0.04515"
Set-TextValue "E23" "-3.07%"
Set-TextValue "D24" "0.001235"
Set-TextValue "E24" "1.00%"
Set-TextValue "D25" "0.004467"
Set-TextValue "E25" "3.07%"
Set-TextValue "D26" "0.0001400"
Set-TextValue "E26" "7.64%"
Set-TextValue "E27" "39.54%"
Set-TextValue "D39" "0.01838"
Set-TextValue "E39" "5.42%"
Set-TextValue "D40" "0.04751"
Set-TextValue "E40" "-0.02%"
Set-TextValue "D41" "0.008123"
Set-TextValue "E41" "4.75%"
Set-TextValue "D42" "0.1396"
Set-TextValue "E42" "-1.05%"
Set-TextValue "D43" "0.007553"
Set-TextValue "E43" "-1.26%"
Set-TextValue "E44" "-3.09%"
Set-TextValue "D45" "0.01047"
Set-TextValue "E45" "4.91%"
Set-TextValue "D46" "0.00006283"
Set-TextValue "E46" "3.27%"
Set-TextValue "E47" "-0.16%"
Set-TextValue "D48" "0.0005802"
Set-TextValue "E48" "0.02%"
Set-TextValue "D49" "30.63"
Set-TextValue "E49" "428.44%"
Set-TextValue "D50" "0.001721"
Set-TextValue "E50" "-36.13%"
Set-TextValue "E51" "-0.16%"
